$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E4").Value = "['Normal']"

$ws.Range("D12").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E12").Value = "['Normal', 'HardwareFault']"

$ws.Range("D15").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E15").Value = "[]"

$ws.Range("D27").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E27").Value = "['SoftwareFault']"

$ws.Range("D28").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E28").Value = "['Normal', 'SoftwareFault']"

$ws.Range("D38").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E38").Value = "['Normal', 'SoftwareFault']"

$ws.Range("D39").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E39").Value = "['Normal', 'SoftwareFault']"

$ws.Range("D40").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E40").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"

$ws.Range("D56").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E56").Value = "[]"

$ws.Range("D73").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal', 'ParamViolation']"

$ws.Range("D88").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E88").Value = "['Normal', 'HardwareFault']"

$ws.Range("D109").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E109").Value = "['Normal', 'SurroundingEnvironment']"

$ws.Range("D113").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E113").Value = "['Normal', 'SoftwareFault']"
